$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 26 with data, mirroring the existing rows' structure.
$row = 26

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"

# Column D holds a date value (serial 44753 = 2022-07-11), styled like the
# existing date cells (same number format as row 25's D cell).
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2022 -Month 7 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(25, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112042
$ws.Cells.Item($row, 7).Value = "Locoto"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 130
$ws.Cells.Item($row, 11).Value = 2700
$ws.Cells.Item($row, 12).Value = 3300
$ws.Cells.Item($row, 13).Value = 2931
$ws.Cells.Item($row, 14).Value = "$/kilo"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 2931
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
